$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.514.59"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.866.00"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.59%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.45%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.91"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.45%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3730"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.22%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07385"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8892"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.55%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07929"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.04"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.865.64"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.71%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.421"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.603"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.76"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.18%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008938"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.75%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.89"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.551.60"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.156"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.57"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.127.11"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.46"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.887"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.54"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.29%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.181"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.95"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08908"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7578"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.83%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.026"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.170"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.490"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.626"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +10.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.082"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01966"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.40%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05281"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.990"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.176"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5226"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.65%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.389"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4880"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.60%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.33"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.89%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.98"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.83%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.656"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06269"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.91"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.29%  "
